$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-18 Saturday", "2025-10-19 Sunday"),
    @("497×2=", "773×4="),
    @("254×4=", "698×4="),
    @("284×6=", "747×5="),
    @("245×7=", "649×5="),
    @("994×5=", "558×5="),
    @("578×3=", "889×3="),
    @("754×8=", "585×4="),
    @("338×9=", "514×2="),
    @("441×8=", "641×7="),
    @("742×5=", "419×7="),
    @("684×8=", "978×4="),
    @("193×9=", "855×4="),
    @("435×9=", "990×6="),
    @("609×8=", "766×6="),
    @("395×7=", "545×4="),
    @("646×5=", "354×5="),
    @("106×2=", "329×9="),
    @("636×2=", "130×2="),
    @("862×6=", "110×3="),
    @("911×2=", "760×5="),
    @("686×3=", "107×6="),
    @("933×4=", "751×4="),
    @("549×3=", "136×7="),
    @("817×4=", "514×3="),
    @("345×9=", "543×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
